# Rady dra Gawkowskiego.docx edit script
#
# Changes applied (per commit: "Add state diagram for Cipher protocol,
# add Cipher interfaces in VisioDiagrams.vsdx"):
#   1. "Dodać FU ..." bullet gets a trailing space run.
#   2. "Nie powinniśmy blokować klienta ..." bullet gets "// " + a green
#      "DONE" marker appended (two runs, DONE colored 00B050).
#   3. "Szyfrowanie ..." bullet gets a single green "// DONE" run appended.
#   4. "Mamy stworzyć architekt uje szyfrowania ..." bullet gets "// " +
#      a green "DONE" marker appended (two runs, DONE colored 00B050).

$d = $word.ActiveDocument

$doneGreen = 5287936   # OLE/BGR encoding of RGB(0x00,0xB0,0x50) -> w:color 00B050

function Append-PlainText {
    param($paraIndex, [string]$text)
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $r.Collapse(0)
    $r.Text = $text
}

function Colorize-TextInParagraph {
    param($paraIndex, [string]$needle)
    $p = $d.Paragraphs($paraIndex)
    $search = $p.Range.Duplicate
    $found = $search.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $search.Font.Color = $doneGreen
    }
}

# 1) "Dodać FU – zmiana hasła przez użytkownika" -> append a trailing space run
Append-PlainText 2 " "

# 2) "Nie powinniśmy blokować klienta ... pamięci. " -> append "// " then green "DONE"
Append-PlainText 7 "// "
Append-PlainText 7 "DONE"
Colorize-TextInParagraph 7 "DONE"

# 3) "Szyfrowanie ... bloki szyfrujące" -> append single green "// DONE" run
Append-PlainText 9 "// DONE"
Colorize-TextInParagraph 9 "// DONE"

# 4) "Mamy stworzyć architekt uje szyfrowania ... XORem i # " -> append "// " then green "DONE"
Append-PlainText 14 "// "
Append-PlainText 14 "DONE"
Colorize-TextInParagraph 14 "DONE"
